$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column-level formatting (mirrors selecting whole columns and applying
#     wrap-text / vertical-centering in the real workbook) ---
$ws.Columns.Item(2).WrapText = $true                 # column B - wrap text
$ws.Columns.Item(5).VerticalAlignment = -4108         # column E - vertical center (xlCenter)
$ws.Columns.Item(6).VerticalAlignment = -4108         # column F - vertical center
$ws.Columns.Item(7).VerticalAlignment = -4108         # column G - vertical center
$ws.Columns.Item(10).VerticalAlignment = -4108        # column J - vertical center
$ws.Columns.Item(11).VerticalAlignment = -4108        # column K - vertical center

# --- Add the new bug report row (row 4) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Top of section cut off on navigation"
$ws.Range("C4").Value = "Clicking on a navigation link will correctly jump to the right section, but the top with the heading is cut off. "
$ws.Range("D4").Value = "1. Open menu(if on mobile)`n2. Click on a navigation link (testimonials, why adopt? Or contact)"
$ws.Range("E4").Value = "Low"
$ws.Range("F4").Value = "Low"
$ws.Range("G4").Value = "Open"
$ws.Range("I4").Value = "N/A"
$ws.Range("J4").Value = 45405
$ws.Range("K4").Value = "Initial report"

# --- Simplify the "Environment" text on the two existing bugs, then reuse
#     the same string for the new row ---
$ws.Range("H2").Value = "Windows 10, Chrome"
$ws.Range("H3").Value = "Windows 10, Chrome"
$ws.Range("H4").Value = "Windows 10, Chrome"

# Formatting for the new row - mirror the other rows' per-cell alignment.
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true
$ws.Range("D4").WrapText = $true
$ws.Range("E4").VerticalAlignment = -4108
$ws.Range("E4").WrapText = $true
$ws.Range("F4").VerticalAlignment = -4108
$ws.Range("G4").VerticalAlignment = -4108
$ws.Range("H4").VerticalAlignment = -4108
$ws.Range("I4").VerticalAlignment = -4108
$ws.Range("J4").VerticalAlignment = -4108
$ws.Range("K4").VerticalAlignment = -4108

$ws.Rows.Item(4).RowHeight = 57.6

# --- Selection moves to J4 (matches the authored file's saved cursor) ---
$ws.Range("J4").Select()
